$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.05832799636040181
$ws.Range("H2").Value = -9.920386237709192
$ws.Range("I2").Value = -121.4873529270648
$ws.Range("G3").Value = 0.06572805537043631
$ws.Range("H3").Value = -3.913476083281062
$ws.Range("G4").Value = -0.04254921411810851
$ws.Range("H4").Value = 5.220906009204787
$ws.Range("G5").Value = -0.03527212521511149
$ws.Range("H5").Value = -30.29081851642291
$ws.Range("G6").Value = -0.1027906970381103
$ws.Range("H6").Value = 3.065319275806043
$ws.Range("G7").Value = -0.0921511138629255
$ws.Range("H7").Value = -0.8536018232424556
$ws.Range("G8").Value = -0.3623223968642763
$ws.Range("H8").Value = 1.248862358372528
$ws.Range("G9").Value = -0.3738361734413548
$ws.Range("H9").Value = 4.166864336001049
$ws.Range("G10").Value = 0.03614983150458126
$ws.Range("H10").Value = 79.03935028294283
$ws.Range("G11").Value = 0.03150633298148364
$ws.Range("H11").Value = 38.82820171166224
$ws.Range("G12").Value = 0.2045720368979911
$ws.Range("H12").Value = -7.745525931469602
$ws.Range("G13").Value = 0.2230290841146326
$ws.Range("H13").Value = -0.9669855660637973
$ws.Range("G14").Value = -0.04685661568271722
$ws.Range("H14").Value = -11.28488631040816
$ws.Range("G15").Value = -0.05003791853063617
$ws.Range("H15").Value = -4.909442074780701
$ws.Range("G16").Value = 0.2140586096099364
$ws.Range("H16").Value = 0.6986814183315547
$ws.Range("G17").Value = 0.2111149330516633
$ws.Range("H17").Value = -4.277088376294032
$ws.Range("G18").Value = 0.07855223677483977
$ws.Range("H18").Value = 7.57597638487549
$ws.Range("G19").Value = 0.08774828869033485
$ws.Range("H19").Value = 16.47279433365793
$ws.Range("G20").Value = -0.07324943780849476
$ws.Range("H20").Value = 2.3106063500251
$ws.Range("G21").Value = -0.08983309954337849
$ws.Range("H21").Value = -3.770719906027618
$ws.Range("G22").Value = 0.07559424518997056
$ws.Range("H22").Value = 2.846658992424517
$ws.Range("G23").Value = 0.0751301975600901
$ws.Range("H23").Value = 9.949875941860485
$ws.Range("G24").Value = 0.05697223742975668
$ws.Range("H24").Value = -14.47167893955226
$ws.Range("G25").Value = 0.06682575910596951
$ws.Range("H25").Value = 21.98921741022475
$ws.Range("G26").Value = 0.1123831206435106
$ws.Range("H26").Value = -5.836589816656104
$ws.Range("G27").Value = 0.1195293469803164
$ws.Range("H27").Value = 4.979057227941124
$ws.Range("G28").Value = 0.1370710158776816
$ws.Range("H28").Value = 6.048303758651891
$ws.Range("G29").Value = 0.1499297765495334
$ws.Range("H29").Value = -0.6037086148149298
$ws.Range("G30").Value = 0.08476799248276795
$ws.Range("H30").Value = 0.5464360990660588
$ws.Range("G31").Value = 0.09043271579595241
$ws.Range("H31").Value = 10.70495823095656
$ws.Range("G32").Value = 0.05485031655421222
$ws.Range("H32").Value = 2.793145303894068
$ws.Range("G33").Value = 0.05767473718886911
$ws.Range("H33").Value = 4.402060494797806
$ws.Range("G34").Value = 0.01434228980352922
$ws.Range("H34").Value = -17.37000580799149
$ws.Range("G35").Value = 0.02134725696142655
$ws.Range("H35").Value = 26.30998684301618
$ws.Range("G36").Value = -0.02469601278756728
$ws.Range("H36").Value = 14.97532384452229
$ws.Range("G37").Value = -0.02160077912148109
$ws.Range("H37").Value = 35.06235597745384
$ws.Range("G38").Value = 0.08295617624488258
$ws.Range("H38").Value = 5.969967666650916
$ws.Range("G39").Value = 0.0678148231732276
$ws.Range("H39").Value = -12.7728999050817
$ws.Range("G40").Value = 0.062747854665526
$ws.Range("H40").Value = -5.227650785716699
$ws.Range("G41").Value = 0.0748304059383727
$ws.Range("H41").Value = 15.08326391983544
$ws.Range("G42").Value = 0.08128983963579929
$ws.Range("H42").Value = 4.497349725639321
$ws.Range("G43").Value = 0.0768170241621656
$ws.Range("H43").Value = -4.175748886807796
$ws.Range("G44").Value = 0.07963406963621836
$ws.Range("H44").Value = -9.760335163190335
$ws.Range("G45").Value = 0.08989991083790649
$ws.Range("H45").Value = -0.5350476005637165
$ws.Range("G46").Value = -0.003000334788855687
$ws.Range("H46").Value = -9.652484789317816
$ws.Range("G47").Value = -0.003372651366724857
$ws.Range("H47").Value = -6959.685158409368
$ws.Range("G48").Value = -0.08985406654980709
$ws.Range("H48").Value = 6.50441530340676
$ws.Range("G49").Value = -0.1045733632872352
$ws.Range("H49").Value = 4.565626178551462
$ws.Range("G50").Value = 0.1758178917829803
$ws.Range("H50").Value = 3.119426953810235
$ws.Range("G51").Value = 0.1768501245636154
$ws.Range("H51").Value = 4.130633477745408
$ws.Range("G52").Value = 0.06213239677001493
$ws.Range("H52").Value = -12.44257797298703
$ws.Range("G53").Value = 0.06798800538368978
$ws.Range("H53").Value = 5.709140987251993
$ws.Range("G54").Value = -0.1333500681251131
$ws.Range("H54").Value = -4.340444856564431
$ws.Range("G55").Value = -0.1375015681986004
$ws.Range("H55").Value = -18.05631474763243
$ws.Range("G56").Value = 0.1870830823541221
$ws.Range("H56").Value = -1.546016135011964
$ws.Range("G57").Value = 0.2034892563218849
$ws.Range("H57").Value = 2.303062542747975
